$d = $word.ActiveDocument

# 1. Header cell: "pvalues" -> "p"
$d.Content.Find.Execute("pvalues", $true, $false, $false, $false, $false,
                         $true, 1, $false, "p", 2)

# 2. p-value for (Constant) row: .205 -> .386
$d.Content.Find.Execute(".205", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".386", 2)

# 3. p-value for var1 row: .119 -> .386
$d.Content.Find.Execute(".119", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".386", 2)

# 4. p-value for var2 row: .232 -> .386
$d.Content.Find.Execute(".232", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".386", 2)

# 5. p-value for var3 row: .480 -> .599
$d.Content.Find.Execute(".480", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".599", 2)

# 6. Add a new paragraph after "Dependent Variable: var5" with the
#    multiple-testing-correction note.
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$newLast = $d.Paragraphs.Last
$newLast.Range.InsertBefore("Multiple tests correction applied to p values: Benjamini-Hochberg")
